$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.452.98"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.49%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.918.15"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.41%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9976"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.81"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9972"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.27%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4700"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.84%  "

$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2880"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +4.83%  "

$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06565"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.62%  "

$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.62"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +11.25%  "

$ws.Range("B11").Value = "Litecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "108.93"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +28.71%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.891.90"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.90%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07564"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.63%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.136"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.43%  "

$ws.Range("B15").Value = "BitcoinCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "317.18"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +29.97%  "

$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.6468"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.43%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.477.19"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.24%  "

$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.97"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.38%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9970"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.30%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007522"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.23%  "

$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.137.35"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.79%  "

$ws.Range("B22").Value = "BinanceUSD"
$ws.Range("C22").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9976"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.17%  "

$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.188"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +5.11%  "

$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.371"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +6.93%  "

$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.277"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.60%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.61"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.61%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.22"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +12.27%  "

$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.030"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +7.65%  "

$ws.Range("B29").Value = "Stellar"
$ws.Range("C29").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1112"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +8.53%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.352"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.55%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.088"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.40%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.942"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.95%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05017"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.40%  "

$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7432"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +5.59%  "

$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.145"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.51%  "

$ws.Range("B36").Value = "Frax"
$ws.Range("C36").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9961"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.23%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.703"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01955"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.28%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.710"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.91%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.012"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8701"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.67%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "107.26"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.95%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.859"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +5.92%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9951"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.45%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4143"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.66%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "68.01"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +8.81%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.262"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.70%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.318"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +8.52%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1200"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.06%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.54"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05627"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.41%  "
